# Update riwayat pembayaran dan logika pembayaran terakhir
#
# - Add a new "Status" column (H) to the Sheet1 table.
# - Mark row 6 (Tiara) as "LUNAS" (paid off) in the new Status column.
# - Row 8 (Dinda)'s Pajak_Terhutang / Tanggal_Jatuh_Tempo / Pajak values were
#   re-entered as plain text (matching the text-typed rows above them)
#   instead of being numeric / date-serial values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New "Status" header in H1, formatted like the other header cells ---
$ws.Range("H1").Value = "Status"
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)  # xlPasteFormats

# --- Row 6 (Tiara) is now fully paid ---
$ws.Range("H6").Value = "LUNAS"

# --- Row 8 (Dinda): store Pajak_Terhutang / Tanggal_Jatuh_Tempo / Pajak as
#     plain text, same as the rest of the "riwayat pembayaran" table, instead
#     of as a number / date serial - and drop the date number-format that
#     used to live on F8. Forcing the Text number format before entering the
#     value stops Excel re-interpreting the digits/date as a number, and
#     clearing formats afterwards removes that temporary format again so the
#     cells end up with plain, unformatted text. ---
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "80000"
$ws.Range("E8").ClearFormats()

$ws.Range("F8").NumberFormat = "@"
$ws.Range("F8").Value = "2026-07-31 00:00:00"
$ws.Range("F8").ClearFormats()

$ws.Range("G8").NumberFormat = "@"
$ws.Range("G8").Value = "80000"
$ws.Range("G8").ClearFormats()
